$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Sheets.Item("ALC")
# Row 28
$ws.Range("H28").Value = 2457
$ws.Range("I28").Value = 3624.75
$ws.Range("K28").Value = 3624.75
$ws.Range("M28").Value = -3139.75
# Row 43
$ws.Range("H43").Value = 4999
$ws.Range("J43").Value = 4998
$ws.Range("L43").Value = 4998
$ws.Range("N43").Value = -5136
# Row 51
$ws.Range("H51").Value = 29999
$ws.Range("J51").Value = 29999
$ws.Range("L51").Value = 29999
$ws.Range("N51").Value = -30967
# Row 97
$ws.Range("H97").Value = 442.42856
$ws.Range("J97").Value = 442.42856
$ws.Range("L97").Value = 1327.28568
$ws.Range("N97").Value = -2319.28568
# Row 118
$ws.Range("H118").Value = 1600
$ws.Range("I118").Value = 1000
$ws.Range("K118").Value = 3000
$ws.Range("M118").Value = -1343
# Row 130
$ws.Range("H130").Value = 88888
$ws.Range("J130").Value = 88888
$ws.Range("L130").Value = 88888
$ws.Range("N130").Value = -98928
# Row 132
$ws.Range("H132").Value = 1773
$ws.Range("I132").Value = 1773
$ws.Range("K132").Value = 5319
$ws.Range("M132").Value = -2789
# Row 138
$ws.Range("H138").Value = 8115.2954
$ws.Range("I138").Value = 1453.8
$ws.Range("K138").Value = 4361.4
$ws.Range("M138").Value = 778.6000000000004

# --- Sheet: ARM ---
$ws = $wb.Sheets.Item("ARM")
# Row 2
$ws.Range("H2").Value = 1254.8572
$ws.Range("I2").Value = 696
$ws.Range("K2").Value = 696
$ws.Range("M2").Value = -583
# Row 32
$ws.Range("H32").Value = 5152.143
$ws.Range("I32").Value = 4733.6665
$ws.Range("K32").Value = 4733.6665
$ws.Range("M32").Value = -4446.6665
# Row 61
$ws.Range("H61").Value = 2996
$ws.Range("I61").Value = 2995.5
$ws.Range("J61").Value = 2996.5
$ws.Range("K61").Value = 2995.5
$ws.Range("L61").Value = 2996.5
$ws.Range("M61").Value = -2783.5
$ws.Range("N61").Value = -3420.5
# Row 74
$ws.Range("H74").Value = 1300.6
$ws.Range("I74").Value = 1332.25
$ws.Range("K74").Value = 1332.25
$ws.Range("M74").Value = -458.25
# Row 77
$ws.Range("H77").Value = 1300.6
$ws.Range("I77").Value = 1332.25
$ws.Range("K77").Value = 6661.25
$ws.Range("M77").Value = -2293.25
# Row 97
$ws.Range("H97").Value = 1358.8889
$ws.Range("I97").Value = 297.83334
$ws.Range("K97").Value = 297.83334
$ws.Range("M97").Value = 198.16666
# Row 116
$ws.Range("H116").Value = 1254.8572
$ws.Range("I116").Value = 696
$ws.Range("K116").Value = 696
$ws.Range("M116").Value = 1598
# Row 122
$ws.Range("H122").Value = 1639
$ws.Range("I122").Value = 1639
$ws.Range("K122").Value = 4917
$ws.Range("M122").Value = -2467
# Row 132
$ws.Range("H132").Value = 1879.871
$ws.Range("I132").Value = 1541.16
$ws.Range("K132").Value = 4623.48
$ws.Range("M132").Value = -2093.48
# Row 136
$ws.Range("H136").Value = 2996
$ws.Range("I136").Value = 2995.5
$ws.Range("J136").Value = 2996.5
$ws.Range("K136").Value = 8986.5
$ws.Range("L136").Value = 8989.5
$ws.Range("M136").Value = -6436.5
$ws.Range("N136").Value = -14089.5

# --- Sheet: BSM ---
$ws = $wb.Sheets.Item("BSM")
# Row 3
$ws.Range("H3").Value = 1254.8572
$ws.Range("I3").Value = 696
$ws.Range("K3").Value = 696
$ws.Range("M3").Value = -582
# Row 20
$ws.Range("H20").Value = 3339.5
$ws.Range("J20").Value = 2487
$ws.Range("L20").Value = 2487
$ws.Range("N20").Value = -2981
# Row 64
$ws.Range("H64").Value = 0
$ws.Range("J64").Value = 0
$ws.Range("L64").Value = 0
$ws.Range("N64").Value = $null
# Row 67
$ws.Range("H67").Value = 0
$ws.Range("J67").Value = 0
$ws.Range("L67").Value = 0
$ws.Range("N67").Value = $null
# Row 99
$ws.Range("H99").Value = 1017.8571
$ws.Range("I99").Value = 854.3333
$ws.Range("K99").Value = 854.3333
$ws.Range("M99").Value = 643.6667
# Row 134
$ws.Range("H134").Value = 2633.5789
$ws.Range("I134").Value = 2235.9412
$ws.Range("K134").Value = 6707.823600000001
$ws.Range("M134").Value = -4172.823600000001

# --- Sheet: CRP ---
$ws = $wb.Sheets.Item("CRP")
# Row 16
$ws.Range("H16").Value = 200
$ws.Range("I16").Value = 200
$ws.Range("J16").Value = 0
$ws.Range("K16").Value = 200
$ws.Range("L16").Value = 0
$ws.Range("M16").Value = 87
$ws.Range("N16").Value = $null
# Row 31
$ws.Range("H31").Value = 1681.0625
$ws.Range("J31").Value = 0
$ws.Range("L31").Value = 0
$ws.Range("N31").Value = $null
# Row 34
$ws.Range("H34").Value = 1681.0625
$ws.Range("J34").Value = 0
$ws.Range("L34").Value = 0
$ws.Range("N34").Value = $null
# Row 58
$ws.Range("H58").Value = 1549.2307
$ws.Range("I58").Value = 1428.4166
$ws.Range("K58").Value = 1428.4166
$ws.Range("M58").Value = -1225.4166
# Row 88
$ws.Range("H88").Value = 57999.332
$ws.Range("J88").Value = 57999.332
$ws.Range("L88").Value = 57999.332
$ws.Range("N88").Value = -58811.332
# Row 91
$ws.Range("H91").Value = 57999.332
$ws.Range("J91").Value = 57999.332
$ws.Range("L91").Value = 57999.332
$ws.Range("N91").Value = -60807.332
# Row 113
$ws.Range("H113").Value = 200
$ws.Range("I113").Value = 200
$ws.Range("J113").Value = 0
$ws.Range("K113").Value = 200
$ws.Range("L113").Value = 0
$ws.Range("M113").Value = 1970
$ws.Range("N113").Value = $null
# Row 132
$ws.Range("H132").Value = 3553.3547
$ws.Range("I132").Value = 3265.074
$ws.Range("K132").Value = 9795.222
$ws.Range("M132").Value = -7265.222
# Row 134
$ws.Range("H134").Value = 4487.7334
$ws.Range("I134").Value = 4562.846
$ws.Range("J134").Value = 3999.5
$ws.Range("K134").Value = 13688.538
$ws.Range("L134").Value = 11998.5
$ws.Range("M134").Value = -11153.538
$ws.Range("N134").Value = -17068.5
# Row 135
$ws.Range("H135").Value = 90000
$ws.Range("J135").Value = 90000
$ws.Range("L135").Value = 90000
$ws.Range("N135").Value = -100140
# Row 136
$ws.Range("H136").Value = 1549.2307
$ws.Range("I136").Value = 1428.4166
$ws.Range("K136").Value = 4285.2498
$ws.Range("M136").Value = -1735.2498

# --- Sheet: CUL ---
$ws = $wb.Sheets.Item("CUL")
# Row 140
$ws.Range("H140").Value = 3199.25
$ws.Range("I140").Value = 3099.3333
$ws.Range("K140").Value = 9297.999899999999
$ws.Range("M140").Value = -4117.999899999999

# --- Sheet: GSM ---
$ws = $wb.Sheets.Item("GSM")
# Row 113
$ws.Range("H113").Value = 0
$ws.Range("I113").Value = 0
$ws.Range("K113").Value = 0
$ws.Range("M113").Value = $null

# --- Sheet: LTW ---
$ws = $wb.Sheets.Item("LTW")
# Row 16
$ws.Range("H16").Value = 850
$ws.Range("I16").Value = 850
$ws.Range("J16").Value = 0
$ws.Range("K16").Value = 850
$ws.Range("L16").Value = 0
$ws.Range("M16").Value = -680
$ws.Range("N16").Value = $null
# Row 93
$ws.Range("H93").Value = 1842.1111
$ws.Range("I93").Value = 1842.1111
$ws.Range("K93").Value = 1842.1111
$ws.Range("M93").Value = -594.1111000000001
# Row 122
$ws.Range("H122").Value = 2683.2222
$ws.Range("I122").Value = 2693.625
$ws.Range("K122").Value = 8080.875
$ws.Range("M122").Value = -5630.875

# --- Sheet: WVR ---
$ws = $wb.Sheets.Item("WVR")
# Row 136
$ws.Range("H136").Value = 738.5789
$ws.Range("I136").Value = 634.6
$ws.Range("J136").Value = 1128.5
$ws.Range("K136").Value = 1903.8
$ws.Range("L136").Value = 3385.5
$ws.Range("M136").Value = 646.1999999999998
$ws.Range("N136").Value = -8485.5
